$d = $word.ActiveDocument

# Work on the first paragraph
$p1 = $d.Paragraphs(1)

# Update the left indent of the first paragraph (120 -> 225 twips = 11.25 pt)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (top/left/bottom/right) with 5pt space, matching the
# border style used elsewhere in the document.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Replace the placeholder text (and drop the trailing space run that used to
# follow it) in one pass, scoped to the first paragraph only.
$found = $p1.Range.Find.Execute("**ID__AFFARS_pgi_5341_topic_3__ID** ", $true, $false, $false, $false, $false,
                                 $true, 1, $false, "**ID__AFFARS_AF_PGI_5341_102__ID**", 2)
